$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NAME values
$ws.Range("A2").Value = "Erik"
$ws.Range("A3").Value = "Linus"

# Update Result/Score/Total for row 2 (Erik, Hep 200M)
$ws.Range("C2").Value = 20.0
$ws.Range("D2").Value = 1398.0
$ws.Range("E2").Value = 1398.0

# Update Result/Score/Total for row 3 (Linus, Hep 200M)
$ws.Range("C3").Value = 19.0
$ws.Range("D3").Value = 1512.0
$ws.Range("E3").Value = 1512.0

# Remove row 4 entirely (the "sten" / "Hep 800M" record)
$ws.Rows.Item(4).Delete()
